# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-08-24 Sunday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2025-08-25 Monday", 2)

# Update the division-problem answers in the single table. Cells are
# addressed by (row, column) rather than text search because several
# problems/answers repeat verbatim elsewhere in the table.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "37÷3=12, 1"
$tbl.Cell(1, 2).Range.Text = "53÷4=13, 1"
$tbl.Cell(1, 3).Range.Text = "29÷2=14, 1"
$tbl.Cell(1, 4).Range.Text = "43÷7=6, 1"
$tbl.Cell(1, 5).Range.Text = "98÷7=14, 0"

$tbl.Cell(5, 1).Range.Text = "84÷2=42, 0"
$tbl.Cell(5, 2).Range.Text = "65÷2=32, 1"
$tbl.Cell(5, 3).Range.Text = "26÷9=2, 8"
$tbl.Cell(5, 4).Range.Text = "41÷9=4, 5"
$tbl.Cell(5, 5).Range.Text = "38÷2=19, 0"

$tbl.Cell(9, 1).Range.Text = "12÷2=6, 0"
$tbl.Cell(9, 2).Range.Text = "95÷8=11, 7"
$tbl.Cell(9, 3).Range.Text = "63÷6=10, 3"
$tbl.Cell(9, 4).Range.Text = "62÷6=10, 2"
$tbl.Cell(9, 5).Range.Text = "35÷4=8, 3"

$tbl.Cell(13, 1).Range.Text = "65÷3=21, 2"
$tbl.Cell(13, 2).Range.Text = "56÷2=28, 0"
$tbl.Cell(13, 3).Range.Text = "50÷5=10, 0"
$tbl.Cell(13, 4).Range.Text = "17÷8=2, 1"
$tbl.Cell(13, 5).Range.Text = "95÷2=47, 1"

$tbl.Cell(17, 1).Range.Text = "66÷8=8, 2"
$tbl.Cell(17, 2).Range.Text = "18÷8=2, 2"
$tbl.Cell(17, 3).Range.Text = "72÷5=14, 2"
$tbl.Cell(17, 4).Range.Text = "70÷2=35, 0"
$tbl.Cell(17, 5).Range.Text = "58÷7=8, 2"
